# Updates cryptos list prices (D) and 1h volume-change percentages (E)
# for the rows whose source data changed, matching the upstream commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.170.88"
$ws.Range("E2").Value = "  +1.34%  "

$ws.Range("D3").Value = "1.604.01"
$ws.Range("E3").Value = "  +0.17%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "'212.09"
$ws.Range("E5").Value = "  +1.46%  "

$ws.Range("E6").Value = "  -0.03%  "

$ws.Range("E7").Value = "  +0.20%  "

$ws.Range("E8").Value = "  +1.04%  "

$ws.Range("D9").Value = "'0.0618"

$ws.Range("D10").Value = "'18.09"
$ws.Range("E10").Value = "  +1.07%  "

$ws.Range("D11").Value = "'0.0794"
$ws.Range("E11").Value = "  +1.43%  "

$ws.Range("D12").Value = "1.826.67"
$ws.Range("E12").Value = "  +0.19%  "

$ws.Range("D13").Value = "1.602.86"
$ws.Range("E13").Value = "  +0.07%  "

$ws.Range("E14").Value = "  -1.29%  "

$ws.Range("E15").Value = "  -0.27%  "

$ws.Range("D16").Value = "26.163.17"
$ws.Range("E16").Value = "  +1.44%  "

$ws.Range("E17").Value = "  -0.02%  "

$ws.Range("E18").Value = "  +1.32%  "

$ws.Range("E19").Value = "  -0.11%  "

$ws.Range("D20").Value = "'198.89"
$ws.Range("E20").Value = "  +4.93%  "

$ws.Range("D21").Value = "'4.23"
$ws.Range("E21").Value = "  +1.12%  "

$ws.Range("D22").Value = "'9.41"
$ws.Range("E22").Value = "  +0.63%  "

$ws.Range("E23").Value = "  +0.91%  "

$ws.Range("E24").Value = "  +2.38%  "

$ws.Range("D25").Value = "'142.06"
$ws.Range("E25").Value = "  +0.71%  "

$ws.Range("E26").Value = "  +3.00%  "

$ws.Range("E27").Value = "  -0.17%  "

$ws.Range("E28").Value = "  +1.18%  "

$ws.Range("E29").Value = "  -0.62%  "

$ws.Range("E30").Value = "  -0.98%  "

$ws.Range("D31").Value = "'0.0472"
$ws.Range("E31").Value = "  +1.00%  "

$ws.Range("E32").Value = "  +1.79%  "

$ws.Range("E33").Value = "  +0.26%  "

$ws.Range("E34").Value = "  +1.61%  "

$ws.Range("E35").Value = "  -1.49%  "

$ws.Range("D36").Value = "1.107.84"
$ws.Range("E36").Value = "  +1.08%  "

$ws.Range("E37").Value = "  -0.16%  "

$ws.Range("E38").Value = "  -0.93%  "

$ws.Range("E39").Value = "  +0.26%  "

$ws.Range("E40").Value = "  +0.72%  "

$ws.Range("D41").Value = "'0.785"
$ws.Range("E41").Value = "  -0.89%  "

$ws.Range("D42").Value = "'0.775"
$ws.Range("E42").Value = "  +4.17%  "

$ws.Range("D43").Value = "1.739.74"
$ws.Range("E43").Value = "  +0.25%  "

$ws.Range("D44").Value = "'5.12"
$ws.Range("E44").Value = "  +0.96%  "

$ws.Range("D45").Value = "'92.82"
$ws.Range("E45").Value = "  -3.03%  "

$ws.Range("D46").Value = "'1.54"
$ws.Range("E46").Value = "  +7.64%  "

$ws.Range("E47").Value = "  -7.59%  "

$ws.Range("D48").Value = "'53.55"
$ws.Range("E48").Value = "  +0.54%  "

$ws.Range("E50").Value = "  -0.06%  "

$ws.Range("E51").Value = "  -0.10%  "
